{"js": "// Reorders the course-requisite entries listed under the \"Requisitos\"\n// heading. The paragraph holding the requisites is a single ListBullet\n// paragraph containing one run per requisite (run text + a manual\n// line break <w:br/>). The content (21 requisites) stays the same,\n// only the order of the runs changes, per the target diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"Requisitos\" heading, the requisites themselves live in\n// the very next paragraph.\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Requisitos\") {\n    headingIndex = i;\n    break;\n  }\n}\nif (headingIndex === -1) {\n  throw new Error('Could not find the \"Requisitos\" heading paragraph.');\n}\n\nconst requisitesParagraph = paragraphs.items[headingIndex + 1];\n\n// Final order of the requisite entries (same 21 items, reordered).\nconst newOrder = [\n  \"LOB1045 -  Leitura e Produ\u00e7\u00e3o de Textos Acad\u00eamicos  (Requisito)\",\n  \"LOB1056 -  Introdu\u00e7\u00e3o aos M\u00e9todos Num\u00e9ricos e Computacionais  (Requisito)\",\n  \"LOQ4095 -  Qu\u00edmica Geral Experimental  (Requisito)\",\n  \"LOQ4098 -  Fundamentos de Qu\u00edmica para Engenharia II (Requisito)\",\n  \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito)\",\n  \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n  \"LOB1037 -  \u00c0lgebra Linear  (Requisito)\",\n  \"LOB1040 -  Laborat\u00f3rio de Eletricidade  (Requisito)\",\n  \"LOB1053 -  F\u00edsica III  (Requisito)\",\n  \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n  \"LOB1009 -  Leitura e Interpreta\u00e7\u00e3o de Desenho T\u00e9cnico  (Requisito)\",\n  \"LOB1011 -  Eletricidade Aplicada  (Requisito)\",\n  \"LOB1018 -  F\u00edsica I  (Requisito)\",\n  \"LOB1024 -  Mec\u00e2nica  (Requisito)\",\n  \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n  \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n  \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n  \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n  \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n  \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n  \"LOB1019 -  F\u00edsica II  (Requisito)\"\n];\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst runsXml = newOrder\n  .map((text) => `<w:r><w:t>${escapeXml(text)}</w:t><w:br/></w:r>`)\n  .join(\"\");\n\nconst ooxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData>' +\n  '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' + runsXml + '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nrequisitesParagraph.getRange().insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Reorders the course-requisite entries listed under the \"Requisitos\"\n# heading. The requisites live in a single ListBullet paragraph that\n# contains one run per requisite (run text + a manual line break\n# <w:br/>). The set of 21 requisites is unchanged - only their order\n# changes, per the target diff.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Requisitos\" heading; the requisites paragraph is the\n# very next paragraph in the document.\n$count = $d.Paragraphs.Count\n$headingIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text\n    $txt = $txt.TrimEnd([char]13, [char]7)\n    if ($txt -eq \"Requisitos\") {\n        $headingIndex = $i\n        break\n    }\n}\nif ($headingIndex -eq -1) {\n    throw \"Could not find the 'Requisitos' heading paragraph.\"\n}\n\n$requisitesIndex = $headingIndex + 1\n\n# Final order of the requisite entries (same 21 items, reordered).\n$newOrder = @(\n    \"LOB1045 -  Leitura e Produ\u00e7\u00e3o de Textos Acad\u00eamicos  (Requisito)\",\n    \"LOB1056 -  Introdu\u00e7\u00e3o aos M\u00e9todos Num\u00e9ricos e Computacionais  (Requisito)\",\n    \"LOQ4095 -  Qu\u00edmica Geral Experimental  (Requisito)\",\n    \"LOQ4098 -  Fundamentos de Qu\u00edmica para Engenharia II (Requisito)\",\n    \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito)\",\n    \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n    \"LOB1037 -  \u00c0lgebra Linear  (Requisito)\",\n    \"LOB1040 -  Laborat\u00f3rio de Eletricidade  (Requisito)\",\n    \"LOB1053 -  F\u00edsica III  (Requisito)\",\n    \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n    \"LOB1009 -  Leitura e Interpreta\u00e7\u00e3o de Desenho T\u00e9cnico  (Requisito)\",\n    \"LOB1011 -  Eletricidade Aplicada  (Requisito)\",\n    \"LOB1018 -  F\u00edsica I  (Requisito)\",\n    \"LOB1024 -  Mec\u00e2nica  (Requisito)\",\n    \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n    \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n    \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n    \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n    \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n    \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n    \"LOB1019 -  F\u00edsica II  (Requisito)\"\n)\n\n$runsXml = \"\"\nforeach ($item in $newOrder) {\n    $escaped = $item.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n    $runsXml += \"<w:r><w:t>\" + $escaped + \"</w:t><w:br/></w:r>\"\n}\n\n$ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData>' +\n    '<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p><w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' + $runsXml + '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$requisitesParagraph = $d.Paragraphs.Item($requisitesIndex)\n$targetRange = $requisitesParagraph.Range\n$null = $targetRange.InsertXML($ooxml)\n\n# InsertXML replaces the range's contents but splits the trailing\n# paragraph mark into its own (now empty) paragraph; merge it back\n# out so the paragraph count matches the original document.\n$firstParagraph = $d.Paragraphs.Item($requisitesIndex)\n$nextParagraph = $d.Paragraphs.Item($requisitesIndex + 1)\n$mergeRange = $d.Range($firstParagraph.Range.End - 1, $nextParagraph.Range.End)\n$null = $mergeRange.Delete()\n"}
